$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2, column C: sample label tweak "Blank3" -> "ABC001"
$ws.Range("C2").Value = "ABC001"

# New custom width for column C (widened to fit the new label)
$ws.Columns("C").ColumnWidth = 12.17

# Selection / zoom state left behind by the editing session
$ws.Range("E9:E10").Select() | Out-Null
$excel.ActiveWindow.Zoom = 169
